# Applies:
#  1) Date placeholder cached text 2017-03-01 -> 2017-03-13 (Master + Layout)
#  2) Reposition/resize a batch of shapes on the "Pillars" slide
#  3) Add a new "REFLECTION" pillar rectangle (4th pillar) after "Rectangle 5"
#
# EMU <-> point helper. The host stores Shape.Left/Top/Width/Height as
# single-precision points and truncates when converting back to EMU, so a
# plain emu/12700.0 can land 1 EMU short after the float32 round-trip. A
# nudge of half an EMU (in point-space) keeps the truncated result exact
# for realistic slide-sized coordinates.
function EmuToPt($emu) {
    return ($emu + 0.5) / 12700.0
}

function SetShapePos($shape, $x, $y) {
    $shape.Left = EmuToPt($x)
    $shape.Top = EmuToPt($y)
}

function SetShapeRect($shape, $x, $y, $cx, $cy) {
    $shape.Left = EmuToPt($x)
    $shape.Top = EmuToPt($y)
    $shape.Width = EmuToPt($cx)
    $shape.Height = EmuToPt($cy)
}

$p = $ppt.ActivePresentation

# --- 1. Refresh the cached "datetimeFigureOut" placeholder text -----------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $msh = $master.Shapes.Item($i)
    if ($msh.Name -eq "Date Placeholder 3") {
        $msh.TextFrame.TextRange.Text = "2017-03-13"
    }
}

$layout = $master.CustomLayouts.Item(1)
for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
    $lsh = $layout.Shapes.Item($i)
    if ($lsh.Name -eq "Date Placeholder 2") {
        $lsh.TextFrame.TextRange.Text = "2017-03-13"
    }
}

# --- 2. Move / resize shapes on the "Pillars" slide ------------------------
$slide = $p.Slides.Item(2)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    switch ($sh.Id) {
        12 { SetShapePos $sh 1275080 783108 }                  # Graphic 11
        14 { SetShapePos $sh 1942464 1437640 }                 # Graphic 13
        15 { SetShapePos $sh 601344 1437640 }                  # Graphic 14
        16 { SetShapePos $sh 1219200 1394460 }                 # Oval 15
        18 { SetShapePos $sh 7592060 792392 }                  # Graphic 17
        19 { SetShapePos $sh 8259444 1446924 }                 # Graphic 18
        20 { SetShapePos $sh 6918324 1446924 }                 # Graphic 19
        21 { SetShapePos $sh 7536180 1403744 }                 # Oval 20
        3  { SetShapeRect $sh 660400 1884680 8519160 690880 }  # Rectangle 2 "SELF ORGANISED TEAMS"
        7  { SetShapeRect $sh 660400 4211320 8519160 690880 }  # Rectangle 6 "COMMON INFRASTRUCTURE"
        4  { SetShapePos $sh 802640 2395220 }                  # Rectangle 3 "VISION/MISSION/GOAL"
        5  { SetShapePos $sh 5120640 2395220 }                 # Rectangle 4 "AUTHORITY/AUTONOMY"
        6  { SetShapePos $sh 2910840 2395220 }                 # Rectangle 5 "FRAMEWORK/WITH CLEAR BOUNDARIES"
    }
}

# --- 3. Add the new "REFLECTION" pillar -------------------------------------
# Duplicate "Rectangle 5" (same line/fill/style refs) via the simulated
# clipboard, then reposition and retext it as the new 4th pillar.
# Slide-scoped shape ids are handed out from a running counter that skips
# ids already present on the slide (1,2,3,4,5,6,7,12,14,15,16,18,19,20,21
# are taken) - churn through 5 throwaway shapes first so the real paste
# lands on id 17, matching the target deck.
for ($k = 1; $k -le 5; $k++) {
    $tmp = $slide.Shapes.AddShape(1, 0, 0, 10, 10)
    $tmp.Delete()
}

$template = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.Id -eq 6) { $template = $sh }
}
$template.Copy()
$pastedRange = $slide.Shapes.Paste()
$newShape = $pastedRange.Item(1)
$newShape.Name = "Rectangle 16"

SetShapeRect $newShape 7330440 2409495 1706880 1998980

# The template has two paragraphs ("FRAMEWORK" / "WITH CLEAR BOUNDARIES").
# Drop the first paragraph (text + its trailing paragraph mark) so the
# remaining paragraph keeps its original endParaRPr, then retext it -
# straight whole-range text assignment would collapse both paragraphs
# into one run and silently drop that endParaRPr.
$newTr = $newShape.TextFrame.TextRange
$firstParaLen = 10  # len("FRAMEWORK") + 1 for the paragraph mark
$newTr.Characters(1, $firstParaLen).Delete()
$newTr.Text = "REFLECTION"
